$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that sat after the Date paragraph.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 2) Wrap the "Location" run (paragraph 5, the centered "Location: Oosterhout"
#    line) with spell-check proof-error markers, matching the diff's
#    <w:proofErr w:type="spellStart"/> ... <w:proofErr w:type="spellEnd"/>.
#    There is no direct object-model property for proofErr markers, so the
#    paragraph's exact OOXML is rebuilt (preserving every run/rPr/rsid) and
#    applied with InsertXML over the paragraph's own range.
# ---------------------------------------------------------------------------
$pLocation = $d.Paragraphs(5)
$locationXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00693EC2" w:rsidRPr="00C50B14" w:rsidRDefault="00DA671B"><w:pPr><w:jc w:val="center"/><w:rPr><w:i/><w:iCs/><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00C50B14"><w:rPr><w:i/><w:iCs/><w:lang w:val="nl-NL"/></w:rPr><w:t>Location</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00693EC2" w:rsidRPr="00C50B14"><w:rPr><w:i/><w:iCs/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r w:rsidR="009608EB" w:rsidRPr="00C50B14"><w:rPr><w:i/><w:iCs/><w:lang w:val="nl-NL"/></w:rPr><w:t>Oosterhout</w:t></w:r></w:p>'
$pLocation.Range.InsertXML($locationXml)

# ---------------------------------------------------------------------------
# 3) Fix the "- Suporting concepts" typo to "- Supporting concepts", split
#    across three runs ("- Su" / "p" / "porting concepts") with a fresh
#    "_GoBack" bookmark placed between the 2nd and 3rd run - exactly as shown
#    in the diff.
# ---------------------------------------------------------------------------
$pSupporting = $d.Paragraphs(27)
$supportingXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00336AAC" w:rsidRPr="00573ACC" w:rsidRDefault="00336AAC" w:rsidP="00336AAC"><w:pPr><w:ind w:left="284"/></w:pPr><w:r><w:t>- Su</w:t></w:r><w:r><w:t>p</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>porting concepts</w:t></w:r></w:p>'
$pSupporting.Range.InsertXML($supportingXml)

# ---------------------------------------------------------------------------
# 4) Delete the entire "Project Elements" / "Plan" / "Legend" section (every
#    paragraph from "Project Elements" through the final "Unorganized" line),
#    leaving only the "- Supporting concepts" paragraph before the sectPr.
# ---------------------------------------------------------------------------
$legendStart = $d.Paragraphs(28)
$legendEnd = $d.Paragraphs($d.Paragraphs.Count)
$legendRange = $d.Range($legendStart.Range.Start, $legendEnd.Range.End)
$legendRange.Delete()
